$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "B"
$ws.Range("B4").Value = "C"
$ws.Range("B5").Value = "D"
$ws.Range("B7").Value = "B"
$ws.Range("B8").Value = "C"
$ws.Range("B9").Value = "D"
